$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data block: rows 16-36 hold 3 workers x 7 periods (2107 down to 2101)
# Worker 1: PEDRO ANTONIO PATERNINA CORONADO (CC 73350956) - Valor Mora 738000
# Worker 2: ORLANDO CABARCAS GUZMAN (CC 9076063) - Valor Mora 877803
# Worker 3: WINDER PACHECO RIVERO (CC 73189894) - Valor Mora 877803
# Salario Basico is 29260 for period 2107, 35112 for all other periods

$workers = @(
    @{ Doc = "73350956"; Name = "PEDRO ANTONIO PATERNINA CORONADO"; Mora = 738000 },
    @{ Doc = "9076063";  Name = "ORLANDO CABARCAS GUZMAN";          Mora = 877803 },
    @{ Doc = "73189894"; Name = "WINDER PACHECO RIVERO";            Mora = 877803 }
)

$periods = @("2107","2106","2105","2104","2103","2102","2101")

$row = 16
foreach ($worker in $workers) {
    foreach ($period in $periods) {
        if ($period -eq "2107") {
            $salario = 29260
        } else {
            $salario = 35112
        }

        $ws.Range("C$row").Value = $worker.Doc
        $ws.Range("D$row").Value = $worker.Name
        $ws.Range("E$row").Value = $period
        $ws.Range("F$row").Value = $salario
        $ws.Range("G$row").Value = $worker.Mora

        $row = $row + 1
    }
}
